$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.307.59'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.589.99'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.41'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '1.813.91'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.08'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.30%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.588.77'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.64'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').Value = '26.306.17'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.53'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '212.24'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('E23').Value = '  +2.27%  '
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.01'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.06'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('D34').Value = '1.327.27'
$ws.Range('E34').Value = '  +3.46%  '
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.604'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.71'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.96%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.995'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -24.19%  '
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.766'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.94'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.726.96'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.13'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  +8.12%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.48'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.86%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0504'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0978'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.49%  '
